$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 8087.6
$ws.Range("I18").Value = 8087.6
$ws.Range("K18").Value = 8087.6
$ws.Range("M18").Value = -7803.6

$ws.Range("H28").Value = 5008.647
$ws.Range("J28").Value = 3189.5715
$ws.Range("L28").Value = 3189.5715
$ws.Range("N28").Value = -4159.5715

$ws.Range("H33").Value = 1350.4
$ws.Range("I33").Value = 786.2857
$ws.Range("K33").Value = 786.2857
$ws.Range("M33").Value = -557.2857

$ws.Range("H40").Value = 3477.9048
$ws.Range("I40").Value = 3008.5625
$ws.Range("J40").Value = 4979.8
$ws.Range("K40").Value = 3008.5625
$ws.Range("L40").Value = 4979.8
$ws.Range("M40").Value = -2833.5625
$ws.Range("N40").Value = -5329.8

$ws.Range("H76").Value = 5931.357
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 5931.357
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 5931.357
$ws.Range("M76").Value = ""
$ws.Range("N76").Value = -6561.357

$ws.Range("H79").Value = 5931.357
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 5931.357
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 5931.357
$ws.Range("M79").Value = ""
$ws.Range("N79").Value = -8115.357

$ws.Range("H82").Value = 1000
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2594

$ws.Range("H85").Value = 1000
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1596

$ws.Range("H94").Value = 55695190
$ws.Range("I94").Value = 71436670
$ws.Range("K94").Value = 71436670
$ws.Range("M94").Value = -71436219

$ws.Range("H100").Value = 94014.75
$ws.Range("I100").Value = 125781.25
$ws.Range("J100").Value = 62248.25
$ws.Range("K100").Value = 125781.25
$ws.Range("L100").Value = 62248.25
$ws.Range("M100").Value = -125240.25
$ws.Range("N100").Value = -63330.25

$ws.Range("H112").Value = 2456.9092
$ws.Range("I112").Value = 1340.25
$ws.Range("J112").Value = 3095
$ws.Range("K112").Value = 4020.75
$ws.Range("L112").Value = 9285
$ws.Range("M112").Value = -2912.75
$ws.Range("N112").Value = -11501

$ws.Range("I113").Value = 16317
$ws.Range("J113").Value = 6629.909
$ws.Range("K113").Value = 16317
$ws.Range("L113").Value = 6629.909
$ws.Range("M113").Value = -13063
$ws.Range("N113").Value = -13137.909

$ws.Range("H125").Value = 6821
$ws.Range("I125").Value = 9062.799999999999
$ws.Range("K125").Value = 81565.2
$ws.Range("M125").Value = -79105.2

$ws.Range("H133").Value = 94949.5
$ws.Range("J133").Value = 94949.5
$ws.Range("L133").Value = 94949.5
$ws.Range("N133").Value = -105069.5

$ws.Range("H137").Value = 10520.833
$ws.Range("I137").Value = 13157.223
$ws.Range("K137").Value = 39471.669
$ws.Range("M137").Value = -36921.669

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2318.1667
$ws.Range("I74").Value = 1251.8572
$ws.Range("J74").Value = 6050.25
$ws.Range("K74").Value = 1251.8572
$ws.Range("L74").Value = 6050.25
$ws.Range("M74").Value = -377.8571999999999
$ws.Range("N74").Value = -7798.25

$ws.Range("H77").Value = 2318.1667
$ws.Range("I77").Value = 1251.8572
$ws.Range("J77").Value = 6050.25
$ws.Range("K77").Value = 6259.286
$ws.Range("L77").Value = 30251.25
$ws.Range("M77").Value = -1891.286
$ws.Range("N77").Value = -38987.25

$ws.Range("H132").Value = 4466.5137
$ws.Range("I132").Value = 3552.4092
$ws.Range("K132").Value = 10657.2276
$ws.Range("M132").Value = -8127.2276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 14995
$ws.Range("J8").Value = 14995
$ws.Range("L8").Value = 14995
$ws.Range("N8").Value = -15275

$ws.Range("H134").Value = 7809
$ws.Range("I134").Value = 8250.941000000001
$ws.Range("K134").Value = 24752.823
$ws.Range("M134").Value = -22217.823

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7292.9697
$ws.Range("I31").Value = 8696.684999999999
$ws.Range("K31").Value = 8696.684999999999
$ws.Range("M31").Value = -8401.684999999999

$ws.Range("H34").Value = 7292.9697
$ws.Range("I34").Value = 8696.684999999999
$ws.Range("K34").Value = 8696.684999999999
$ws.Range("M34").Value = -8494.684999999999

$ws.Range("H99").Value = 5053601
$ws.Range("I99").Value = 10559189
$ws.Range("K99").Value = 10559189
$ws.Range("M99").Value = -10557691

$ws.Range("H126").Value = 5053601
$ws.Range("I126").Value = 10559189
$ws.Range("K126").Value = 31677567
$ws.Range("M126").Value = -31675097

$ws.Range("H132").Value = 1773.5
$ws.Range("I132").Value = 1653.25
$ws.Range("K132").Value = 4959.75
$ws.Range("M132").Value = -2429.75

$ws.Range("H141").Value = 306932.34
$ws.Range("J141").Value = 326653.47
$ws.Range("L141").Value = 326653.47
$ws.Range("N141").Value = -337013.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").Value = ""

$ws.Range("H132").Value = 28574.111
$ws.Range("I132").Value = 702.3333
$ws.Range("J132").Value = 42510
$ws.Range("K132").Value = 6320.9997
$ws.Range("L132").Value = 382590
$ws.Range("M132").Value = -3790.9997
$ws.Range("N132").Value = -387650

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 26857142
$ws.Range("I7").Value = 34800000
$ws.Range("J7").Value = 7000000
$ws.Range("K7").Value = 34800000
$ws.Range("L7").Value = 7000000
$ws.Range("M7").Value = -34799888
$ws.Range("N7").Value = -7000224

$ws.Range("H8").Value = 26857142
$ws.Range("I8").Value = 34800000
$ws.Range("J8").Value = 7000000
$ws.Range("K8").Value = 34800000
$ws.Range("L8").Value = 7000000
$ws.Range("M8").Value = -34799861
$ws.Range("N8").Value = -7000278

$ws.Range("H26").Value = 16998
$ws.Range("I26").Value = 16597
$ws.Range("J26").Value = 17666.334
$ws.Range("K26").Value = 16597
$ws.Range("L26").Value = 17666.334
$ws.Range("M26").Value = -16317
$ws.Range("N26").Value = -18226.334

$ws.Range("H50").Value = 16998
$ws.Range("I50").Value = 16597
$ws.Range("J50").Value = 17666.334
$ws.Range("K50").Value = 16597
$ws.Range("L50").Value = 17666.334
$ws.Range("M50").Value = -16099
$ws.Range("N50").Value = -18662.334

$ws.Range("H80").Value = 6159.3438
$ws.Range("J80").Value = 7120
$ws.Range("L80").Value = 7120
$ws.Range("N80").Value = -9116

$ws.Range("H83").Value = 6159.3438
$ws.Range("J83").Value = 7120
$ws.Range("L83").Value = 35600
$ws.Range("N83").Value = -45584

$ws.Range("H97").Value = 5128.4287
$ws.Range("I97").Value = 6255.727
$ws.Range("J97").Value = 995
$ws.Range("K97").Value = 6255.727
$ws.Range("L97").Value = 995
$ws.Range("M97").Value = -5759.727
$ws.Range("N97").Value = -1987

$ws.Range("H123").Value = 42258.668
$ws.Range("J123").Value = 42258.668
$ws.Range("L123").Value = 42258.668
$ws.Range("N123").Value = -47158.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 46248.5
$ws.Range("J38").Value = 46248.5
$ws.Range("L38").Value = 46248.5
$ws.Range("N38").Value = -47068.5

$ws.Range("H82").Value = 3411.5833
$ws.Range("I82").Value = 3898.75
$ws.Range("K82").Value = 3898.75
$ws.Range("M82").Value = -3537.75

$ws.Range("H85").Value = 3411.5833
$ws.Range("I85").Value = 3898.75
$ws.Range("K85").Value = 3898.75
$ws.Range("M85").Value = -2650.75

$ws.Range("H132").Value = 713124.7
$ws.Range("J132").Value = 6811.25
$ws.Range("L132").Value = 20433.75
$ws.Range("N132").Value = -25493.75

$ws.Range("H136").Value = 5732.3335
$ws.Range("J136").Value = 7959
$ws.Range("L136").Value = 23877
$ws.Range("N136").Value = -28977

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8205.4375
$ws.Range("I81").Value = 10318
$ws.Range("J81").Value = 3557.8
$ws.Range("K81").Value = 20636
$ws.Range("L81").Value = 7115.6
$ws.Range("M81").Value = -19575
$ws.Range("N81").Value = -9237.6

$ws.Range("H84").Value = 8205.4375
$ws.Range("I84").Value = 10318
$ws.Range("J84").Value = 3557.8
$ws.Range("K84").Value = 103180
$ws.Range("L84").Value = 35578
$ws.Range("M84").Value = -97876
$ws.Range("N84").Value = -46186

$ws.Range("H99").Value = 27998
$ws.Range("I99").Value = 27998
$ws.Range("K99").Value = 27998
$ws.Range("M99").Value = -25003

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = ""

$ws.Range("H132").Value = 13588.84
$ws.Range("I132").Value = 14939.167
$ws.Range("K132").Value = 44817.501
$ws.Range("M132").Value = -42287.501

$ws.Range("H136").Value = 645035.8
$ws.Range("I136").Value = 908416.75
$ws.Range("J136").Value = 5396.4287
$ws.Range("K136").Value = 2725250.25
$ws.Range("L136").Value = 16189.2861
$ws.Range("M136").Value = -2722700.25
$ws.Range("N136").Value = -21289.2861
